# "Generate Report for Handoff"
# b.md has moved from "Handed back: in sync with en-US" to "Ready for handoff"
# for both the zh-cn and de-de locales. Update the Overview sheet plus the two
# per-locale detail sheets accordingly.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19486a4575ff0687249d96f5ce3ec19eef26ea4a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/047f273e17e351d16de38712deae5506f818a7cf/e2e/b.md."
$pasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# Helper: Excel auto-coerces a literal "True"/"False" assigned straight to
# .Value into a real Boolean cell. The source report keeps this column as
# plain text, so round-trip the value through a text formula and paste it
# back as a value to keep the cell type as a string.
function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy($range)
    $range.PasteSpecial($pasteValues)
}

# ---- Overview sheet: row 3 is b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-03 06:40:40"

# ---- zh-cn sheet: row 3 is b.md ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsZhCn.Range("F3") "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-03 06:40:35"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14

# ---- de-de sheet: row 3 is b.md ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsDeDe.Range("F3") "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-03 06:40:40"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14
